$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44160
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 7500
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7667
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 128

# Row 3
$ws.Range("D3").Value = 44253
$ws.Range("J3").Value = 95
$ws.Range("K3").Value = 9500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9658
$ws.Range("P3").Value = 161

# Row 4
$ws.Range("D4").Value = 44266
$ws.Range("J4").Value = 60
$ws.Range("M4").Value = 9208
$ws.Range("P4").Value = 153

# Row 5
$ws.Range("D5").Value = 44159
$ws.Range("J5").Value = 35
$ws.Range("K5").Value = 7500
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 7714
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 129

# Row 6
$ws.Range("D6").Value = 44224
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 8500
$ws.Range("L6").Value = 9000
$ws.Range("M6").Value = 8719
$ws.Range("P6").Value = 145

# Row 7
$ws.Range("D7").Value = 44218
$ws.Range("J7").Value = 65
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 9615
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 160

# Row 9
$ws.Range("D9").Value = 44202
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 9000
$ws.Range("M9").Value = 8400
$ws.Range("P9").Value = 140

# Row 10
$ws.Range("D10").Value = 44264
$ws.Range("J10").Value = 43
$ws.Range("K10").Value = 8500
$ws.Range("M10").Value = 8709
$ws.Range("P10").Value = 145

# Row 11
$ws.Range("D11").Value = 44208
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 7000
$ws.Range("M11").Value = 7350
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 122

# Row 12
$ws.Range("D12").Value = 44204
$ws.Range("J12").Value = 45
$ws.Range("K12").Value = 9500
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 9722
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 162

# Row 13
$ws.Range("D13").Value = 44259
$ws.Range("J13").Value = 70
$ws.Range("K13").Value = 9000
$ws.Range("L13").Value = 9500
$ws.Range("M13").Value = 9214
$ws.Range("P13").Value = 154

# Row 14
$ws.Range("D14").Value = 44216
$ws.Range("J14").Value = 55
$ws.Range("K14").Value = 9500
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9773
$ws.Range("P14").Value = 163

# Row 15
$ws.Range("D15").Value = 44210
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 8417
$ws.Range("O15").Value = "Región de Arica y Parinacota"
$ws.Range("P15").Value = 140

# Row 16
$ws.Range("D16").Value = 44271
$ws.Range("J16").Value = 55
$ws.Range("L16").Value = 9500
$ws.Range("M16").Value = 9227
$ws.Range("P16").Value = 154
